$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, $val) {
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "67.237.66"

Set-TextValue $ws.Range("D3") "3.528.89"
$ws.Range("E3").Value = "  +10.88%  "

$ws.Range("E4").Value = "  -0.04%  "

Set-TextValue $ws.Range("D5") "190.71"
$ws.Range("E5").Value = "  +10.51%  "

Set-TextValue $ws.Range("D6") "553.43"
$ws.Range("E6").Value = "  +4.81%  "

Set-TextValue $ws.Range("D7") "3.524.05"
$ws.Range("E7").Value = "  +10.88%  "

$ws.Range("E8").Value = "  +2.50%  "

$ws.Range("E9").Value = "  -0.01%  "

$ws.Range("E10").Value = "  +4.62%  "

Set-TextValue $ws.Range("D11") "0.153"
$ws.Range("E11").Value = "  +15.91%  "

Set-TextValue $ws.Range("D12") "55.15"
$ws.Range("E12").Value = "  +3.19%  "

$ws.Range("E13").Value = "  +7.62%  "

$ws.Range("E14").Value = "  +3.31%  "

Set-TextValue $ws.Range("D15") "4.086.57"
$ws.Range("E15").Value = "  +10.75%  "

Set-TextValue $ws.Range("D16") "3.530.66"
$ws.Range("E16").Value = "  +11.01%  "

$ws.Range("E17").Value = "  +3.40%  "

Set-TextValue $ws.Range("D18") "67.276.23"
$ws.Range("E18").Value = "  +7.75%  "

$ws.Range("E19").Value = "  +5.83%  "

$ws.Range("E20").Value = "  +7.78%  "

$ws.Range("E21").Value = "  +3.17%  "

Set-TextValue $ws.Range("D22") "432.21"
$ws.Range("E22").Value = "  +18.02%  "

$ws.Range("E23").Value = "  +4.12%  "

Set-TextValue $ws.Range("D24") "85.05"
$ws.Range("E24").Value = "  +4.82%  "

$ws.Range("E25").Value = "  +7.66%  "

Set-TextValue $ws.Range("D26") "11.20"
$ws.Range("E26").Value = "  -0.49%  "

Set-TextValue $ws.Range("D27") "2.91"
$ws.Range("E27").Value = "  +9.97%  "

$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D28") "12.03"
$ws.Range("E28").Value = "  +5.98%  "

$ws.Range("B29").Value = "Filecoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D29") "8.99"
$ws.Range("E29").Value = "  +9.85%  "

$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D30") "30.32"
$ws.Range("E30").Value = "  +7.00%  "

$ws.Range("B31").Value = "Bittensor"
$ws.Range("C31").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws.Range("D31") "649.61"
$ws.Range("E31").Value = "  +1.30%  "

$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D32") "6.70"
$ws.Range("E32").Value = "  +2.99%  "

$ws.Range("B33").Value = "Cosmos"
$ws.Range("C33").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D33") "11.73"
$ws.Range("E33").Value = "  +3.59%  "

$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D34") "0.111"
$ws.Range("E34").Value = "  +5.52%  "

$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D35") "59.50"
$ws.Range("E35").Value = "  +5.03%  "

$ws.Range("B36").Value = "PEPE"
$ws.Range("C36").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue $ws.Range("D36") "0.0₃0826"
$ws.Range("E36").Value = "  +16.87%  "

Set-TextValue $ws.Range("D37") "38.71"
$ws.Range("E37").Value = "  +4.52%  "

$ws.Range("B38").Value = "Dai"
$ws.Range("C38").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Range("D38") "0.999"
$ws.Range("E38").Value = "  -0.15%  "

$ws.Range("B39").Value = "TheGraph"
$ws.Range("C39").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue $ws.Range("D39") "0.392"
$ws.Range("E39").Value = "  +4.32%  "

$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D40") "0.142"
$ws.Range("E40").Value = "  +15.07%  "

$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D41") "3.34"
$ws.Range("E41").Value = "  +14.38%  "

$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws.Range("D42") "0.999"
$ws.Range("E42").Value = "  -0.03%  "

$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue $ws.Range("D43") "3.025.26"
$ws.Range("E43").Value = "  +5.34%  "

$ws.Range("B44").Value = "Fetch.AI"
$ws.Range("C44").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws.Range("D44") "2.66"
$ws.Range("E44").Value = "  +4.72%  "

$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D45") "2.89"
$ws.Range("E45").Value = "  +9.24%  "

$ws.Range("B46").Value = "ThetaToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue $ws.Range("D46") "2.88"
$ws.Range("E46").Value = "  +10.93%  "

$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue $ws.Range("D47") "3.36"
$ws.Range("E47").Value = "  +12.04%  "

$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D48") "0.0418"
$ws.Range("E48").Value = "  +6.60%  "

$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D49") "0.131"
$ws.Range("E49").Value = "  +5.85%  "

$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue $ws.Range("D50") "8.75"
$ws.Range("E50").Value = "  +14.14%  "

$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D51") "142.76"
$ws.Range("E51").Value = "  +6.38%  "
